$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 149.32
$ws.Range("I15").Value = 149.32
$ws.Range("K15").Value = 447.96
$ws.Range("M15").Value = -278.96
# Row 98
$ws.Range("H98").Value = 2027.7273
$ws.Range("I98").Value = 2130.5
$ws.Range("K98").Value = 2130.5
$ws.Range("M98").Value = -632.5
# Row 116
$ws.Range("H116").Value = 4738.769
$ws.Range("I116").Value = 4889.3335
$ws.Range("K116").Value = 4889.3335
$ws.Range("M116").Value = -1447.3335
# Row 122
$ws.Range("H122").Value = 2027.7273
$ws.Range("I122").Value = 2130.5
$ws.Range("K122").Value = 6391.5
$ws.Range("M122").Value = -3941.5
# Row 131
$ws.Range("H131").Value = 1709.8334
$ws.Range("I131").Value = 981.3077
$ws.Range("K131").Value = 2943.9231
$ws.Range("M131").Value = 2096.0769
# Row 138
$ws.Range("H138").Value = 1408.61
$ws.Range("J138").Value = 3319.5518
$ws.Range("L138").Value = 9958.6554
$ws.Range("N138").Value = -20238.6554

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 1800
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1800
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = 1800
$ws.Range("N3").Value = -2030
# Row 74
$ws.Range("H74").Value = 1297.8148
$ws.Range("I74").Value = 1353.5625
$ws.Range("K74").Value = 1353.5625
$ws.Range("M74").Value = -479.5625
# Row 77
$ws.Range("H77").Value = 1297.8148
$ws.Range("I77").Value = 1353.5625
$ws.Range("K77").Value = 6767.8125
$ws.Range("M77").Value = -2399.8125
# Row 132
$ws.Range("H132").Value = 1046.4912
$ws.Range("I132").Value = 695.12195
$ws.Range("J132").Value = 1946.875
$ws.Range("K132").Value = 2085.36585
$ws.Range("L132").Value = 5840.625
$ws.Range("M132").Value = 444.6341499999999
$ws.Range("N132").Value = -10900.625

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2307.28
$ws.Range("I86").Value = 2219.125
$ws.Range("J86").Value = 2464
$ws.Range("K86").Value = 2219.125
$ws.Range("L86").Value = 2464
$ws.Range("M86").Value = -1096.125
$ws.Range("N86").Value = -4710
# Row 89
$ws.Range("H89").Value = 2307.28
$ws.Range("I89").Value = 2219.125
$ws.Range("J89").Value = 2464
$ws.Range("K89").Value = 11095.625
$ws.Range("L89").Value = 12320
$ws.Range("M89").Value = -5479.625
$ws.Range("N89").Value = -23552
# Row 134
$ws.Range("H134").Value = 109254.055
$ws.Range("I134").Value = 4294.4614
$ws.Range("J134").Value = 336666.5
$ws.Range("K134").Value = 12883.3842
$ws.Range("L134").Value = 1009999.5
$ws.Range("M134").Value = -10348.3842
$ws.Range("N134").Value = -1015069.5

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 1347.5
$ws.Range("I3").Value = 1347.5
$ws.Range("K3").Value = 1347.5
$ws.Range("M3").Value = -1234.5
# Row 31
$ws.Range("H31").Value = 2455.9805
$ws.Range("I31").Value = 2479.8918
$ws.Range("J31").Value = 2392.7856
$ws.Range("K31").Value = 2479.8918
$ws.Range("L31").Value = 2392.7856
$ws.Range("M31").Value = -2184.8918
$ws.Range("N31").Value = -2982.7856
# Row 34
$ws.Range("H34").Value = 2455.9805
$ws.Range("I34").Value = 2479.8918
$ws.Range("J34").Value = 2392.7856
$ws.Range("K34").Value = 2479.8918
$ws.Range("L34").Value = 2392.7856
$ws.Range("M34").Value = -2277.8918
$ws.Range("N34").Value = -2796.7856
# Row 132
$ws.Range("H132").Value = 1105
$ws.Range("I132").Value = 866.4545000000001
$ws.Range("J132").Value = 2536.2727
$ws.Range("K132").Value = 2599.3635
$ws.Range("L132").Value = 7608.8181
$ws.Range("M132").Value = -69.36350000000039
$ws.Range("N132").Value = -12668.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 3096143
$ws.Range("I2").Value = 6536080
$ws.Range("J2").Value = 199.4
$ws.Range("K2").Value = 39216480
$ws.Range("L2").Value = 1196.4
$ws.Range("M2").Value = -39216367
$ws.Range("N2").Value = -1422.4
# Row 11
$ws.Range("H11").Value = 40127
$ws.Range("I11").Value = 66745
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 200235
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -200095
$ws.Range("N11").Value = -880
# Row 68
$ws.Range("H68").Value = 818.2727
$ws.Range("I68").Value = 600.6667
$ws.Range("J68").Value = 899.875
$ws.Range("K68").Value = 1802.0001
$ws.Range("L68").Value = 2699.625
$ws.Range("M68").Value = -991.0001
$ws.Range("N68").Value = -4321.625
# Row 71
$ws.Range("H71").Value = 818.2727
$ws.Range("I71").Value = 600.6667
$ws.Range("J71").Value = 899.875
$ws.Range("K71").Value = 5406.0003
$ws.Range("L71").Value = 8098.875
$ws.Range("M71").Value = -1350.0003
$ws.Range("N71").Value = -16210.875
# Row 122
$ws.Range("H122").Value = 910021.8
$ws.Range("J122").Value = 1429778.6
$ws.Range("L122").Value = 12868007.4
$ws.Range("N122").Value = -12872907.4
# Row 131
$ws.Range("H131").Value = 5176639.5
$ws.Range("I131").Value = 50101344
$ws.Range("J131").Value = 12880.115
$ws.Range("K131").Value = 150304032
$ws.Range("L131").Value = 38640.345
$ws.Range("M131").Value = -150298992
$ws.Range("N131").Value = -48720.345

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3412.3809
$ws.Range("I80").Value = 4068.5715
$ws.Range("J80").Value = 2100
$ws.Range("K80").Value = 4068.5715
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -3070.5715
$ws.Range("N80").Value = -4096
# Row 83
$ws.Range("H83").Value = 3412.3809
$ws.Range("I83").Value = 4068.5715
$ws.Range("J83").Value = 2100
$ws.Range("K83").Value = 20342.8575
$ws.Range("L83").Value = 10500
$ws.Range("M83").Value = -15350.8575
$ws.Range("N83").Value = -20484
# Row 132
$ws.Range("H132").Value = 3657.5386
$ws.Range("I132").Value = 3794.2222
$ws.Range("J132").Value = 3350
$ws.Range("K132").Value = 11382.6666
$ws.Range("L132").Value = 10050
$ws.Range("M132").Value = -8852.6666
$ws.Range("N132").Value = -15110

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 644.2941
$ws.Range("I136").Value = 321.08334
$ws.Range("J136").Value = 1420
$ws.Range("K136").Value = 963.2500200000001
$ws.Range("L136").Value = 4260
$ws.Range("M136").Value = 1586.74998
$ws.Range("N136").Value = -9360
